# Apply the "tuning fonts & object widths" edit:
# On every slide's "Welcome to Master the Tidyverse" title textbox, bump the
# font size to 96pt and split the second run (" Master the Tidyverse") into
# two runs (" Master the " / "Tidyverse") so "Tidyverse" can carry its own
# run-level formatting (matches the canonical OOXML: three <a:r> runs, each
# with rPr sz="9600").

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    $target = $null
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $candidate = $s.Shapes.Item($j)
        if ($candidate.Name -eq "Welcome to Master the Tidyverse") {
            $target = $candidate
        }
    }

    if ($target -eq $null) {
        continue
    }

    $tr = $target.TextFrame.TextRange
    $fullText = $tr.Text

    $marker = "Welcome to"
    if ($fullText.Length -lt $marker.Length) {
        continue
    }

    $run1Len = $marker.Length
    $rest = $fullText.Substring($run1Len)
    $tidyMarker = "Tidyverse"
    $tidyIdx = $rest.LastIndexOf($tidyMarker)
    if ($tidyIdx -lt 0) {
        continue
    }

    $run2Len = $tidyIdx
    $run3Len = $rest.Length - $tidyIdx

    $run1 = $tr.Characters(1, $run1Len)
    $run2 = $tr.Characters($run1Len + 1, $run2Len)
    $run3 = $tr.Characters($run1Len + $run2Len + 1, $run3Len)

    $run1.Font.Size = 96
    $run2.Font.Size = 96
    $run3.Font.Size = 96
}
